# "Processed Results - updated measures"
#
# Adds a Std / Relative std summary block (labels in D15:E15, formulas in
# D16:E16) to the Accelerometer sheet, drops the two stale duplicate
# "_xlchart" defined names left over from the chart tooling, and moves the
# sheet selection onto the newly added D15:E16 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leftover duplicate chart-helper defined names (v1.2 / v1.3),
# keeping v1.0 / v1.1.
$wb.Names.Item("_xlchart.v1.3").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()

# New header labels, bold like the other D/E header pairs (Min/Max,
# Q1/Median, IQR/Q3) on this sheet.
$ws.Range("D15").Value = "Std"
$ws.Range("E15").Value = "Relative std"
$ws.Range("D15:E15").Font.Bold = $true

# New summary formulas: standard deviation of the sample, and that std
# expressed as a percentage of the mean (E3).
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "=(D16 / E3) *100"

# Select the newly added block.
$ws.Range("D15:E16").Select()
